$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.964.02'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '3.789.80'
$ws.Range("E3").Value = '  -1.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '701.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("D7").Value = '3.791.68'
$ws.Range("E7").Value = '  -1.29%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.523'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.40'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.476'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").Value = '4.428.23'
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.818.66'
$ws.Range("E16").Value = '  -1.16%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '70.928.90'
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.114'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '512.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.711'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000139'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.55'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.93%  '
$ws.Range("D27").Value = '3.958.25'
$ws.Range("E27").Value = '  -1.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.27'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.33%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '28.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.174'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("B38").Value = 'RenzoRestakedETH'
$ws.Range("C38").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D38").Value = '3.750.31'
$ws.Range("E38").Value = '  -1.38%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.100'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.45%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.42'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.72%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.31%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '169.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.98%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '49.73'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.25%  '
$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000305'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '419.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.36'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.72%  '
